$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "31.181.49"
$c.Style = "Normal"
$ws.Range("E2").Value = "  +2.14%  "

$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "1.939.37"
$c.Style = "Normal"
$ws.Range("E3").Value = "  +0.93%  "

$ws.Range("E4").Value = "  -0.43%  "

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "242.15"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +1.28%  "

$ws.Range("E6").Value = "  -0.40%  "

$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.4795"
$c.Style = "Normal"
$ws.Range("E7").Value = "  -0.02%  "

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.2909"
$c.Style = "Normal"
$ws.Range("E8").Value = "  +1.08%  "

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.06785"
$c.Style = "Normal"
$ws.Range("E9").Value = "  +1.22%  "

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "20.19"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +7.33%  "

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "104.26"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +0.02%  "

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.07838"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +1.11%  "

$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "1.947.71"
$c.Style = "Normal"
$ws.Range("E13").Value = "  +1.20%  "

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "5.301"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +1.07%  "

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "0.6985"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +2.39%  "

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "296.84"
$c.Style = "Normal"
$ws.Range("E16").Value = "  +11.78%  "

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "31.176.02"
$c.Style = "Normal"
$ws.Range("E17").Value = "  +1.97%  "

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "2.203.18"
$c.Style = "Normal"
$ws.Range("E18").Value = "  +1.43%  "

$ws.Range("E19").Value = "  +2.22%  "

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "0.000007608"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +1.10%  "

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "5.571"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +2.51%  "

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "0.9999"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -0.37%  "

$ws.Range("E23").Value = "  -0.47%  "

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "6.421"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +1.08%  "

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "9.561"
$c.Style = "Normal"
$ws.Range("E25").Value = "  -0.85%  "

$ws.Range("E26").Value = "  +3.51%  "

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "19.82"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +4.17%  "

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "2.099"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +0.16%  "

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "1.392"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +0.20%  "

$ws.Range("E30").Value = "  -1.45%  "

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "4.625"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +0.50%  "

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "1.536"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +1.12%  "

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "4.341"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +1.83%  "

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "0.04845"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +1.91%  "

$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "0.7383"
$c.Style = "Normal"
$ws.Range("E35").Value = "  -0.08%  "

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "1.134"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +1.22%  "

$ws.Range("E37").Value = "  +1.66%  "

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.01963"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +1.16%  "

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "6.809"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +7.40%  "

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "2.634"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -0.38%  "

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "76.62"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +1.40%  "

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "2.036"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +1.76%  "

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "0.8726"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +1.23%  "

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "0.4367"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +1.89%  "

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "105.78"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -0.39%  "

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "0.9999"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -0.40%  "

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "1.029.42"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +3.22%  "

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "7.584"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +0.30%  "

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "9.292"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +3.48%  "

$ws.Range("E50").Value = "  +0.55%  "

$ws.Range("E51").Value = "  +0.11%  "
